$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column at M (13), shifting existing M:V -> N:W.
# New column inherits formatting from the column to its left (L), matching
# native Excel "Insert Column" behaviour.
$ws.Columns("M").Insert()

# Fill in the new header (row 7) and placeholder (row 9) cells for the
# newly inserted "Trang thai dai ly" column.
$ws.Range("M7").Value = "Trạng thái đại lý"
$ws.Range("M9").Value = "{{ReportStoreCheckeds.SalesEmployees.Dates.Contents.StoreStatusName}}"

$ws.Range("N18").Select()
